# Auto-generated Excel COM-interop script
# Applies numeric cell updates across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# to match the scheduled-runner data refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (34 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 41668012
$ws.Range("J51").Value = 38462716
$ws.Range("L51").Value = 38462716
$ws.Range("N51").Value = -38463684
$ws.Range("H62").Value = 2397.2632
$ws.Range("I62").Value = 2452.4666
$ws.Range("K62").Value = 2452.4666
$ws.Range("M62").Value = -1828.4666
$ws.Range("H65").Value = 2397.2632
$ws.Range("I65").Value = 2452.4666
$ws.Range("K65").Value = 12262.333
$ws.Range("M65").Value = -9142.333000000001
$ws.Range("H86").Value = 102566904
$ws.Range("I86").Value = 111113704
$ws.Range("K86").Value = 111113704
$ws.Range("M86").Value = -111112581
$ws.Range("H89").Value = 102566904
$ws.Range("I89").Value = 111113704
$ws.Range("K89").Value = 555568520
$ws.Range("M89").Value = -555562904
$ws.Range("H100").Value = 36905850
$ws.Range("I100").Value = 2083825.1
$ws.Range("J100").Value = 83335220
$ws.Range("K100").Value = 2083825.1
$ws.Range("L100").Value = 83335220
$ws.Range("M100").Value = -2083284.1
$ws.Range("N100").Value = -83336302
$ws.Range("H138").Value = 6056.4614
$ws.Range("I138").Value = 15637.363
$ws.Range("J138").Value = 4104.7964
$ws.Range("K138").Value = 46912.089
$ws.Range("L138").Value = 12314.3892
$ws.Range("M138").Value = -41772.089
$ws.Range("N138").Value = -22594.3892

# --- Sheet: ARM (16 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 73999
$ws.Range("J109").Value = 73999
$ws.Range("L109").Value = 73999
$ws.Range("N109").Value = -76773
$ws.Range("H133").Value = 104995
$ws.Range("J133").Value = 104995
$ws.Range("L133").Value = 104995
$ws.Range("N133").Value = -110055
$ws.Range("H134").Value = 99073.22
$ws.Range("J134").Value = 99073.22
$ws.Range("L134").Value = 99073.22
$ws.Range("N134").Value = -109213.22
$ws.Range("H140").Value = 106853.664
$ws.Range("J140").Value = 106853.664
$ws.Range("L140").Value = 106853.664
$ws.Range("N140").Value = -117213.664

# --- Sheet: BSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 5000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5454
$ws.Range("H22").Value = 1919.7
$ws.Range("J22").Value = 2524.25
$ws.Range("L22").Value = 2524.25
$ws.Range("N22").Value = -2870.25
$ws.Range("I134").Value = 3025.1667
$ws.Range("K134").Value = 9075.500100000001
$ws.Range("M134").Value = -6540.500100000001

# --- Sheet: CRP (26 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 45000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 45000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H94").Value = 1745.9
$ws.Range("I94").Value = 1226.4
$ws.Range("J94").Value = 2265.4
$ws.Range("K94").Value = 1226.4
$ws.Range("L94").Value = 2265.4
$ws.Range("M94").Value = -775.4000000000001
$ws.Range("N94").Value = -3167.4
$ws.Range("H122").Value = 2915.3416
$ws.Range("I122").Value = 3075.1316
$ws.Range("J122").Value = 891.3333
$ws.Range("K122").Value = 9225.3948
$ws.Range("L122").Value = 2673.9999
$ws.Range("M122").Value = -6775.3948
$ws.Range("N122").Value = -7573.9999
$ws.Range("H141").Value = 643642.9399999999
$ws.Range("J141").Value = 741968.5
$ws.Range("L141").Value = 741968.5
$ws.Range("N141").Value = -752328.5

# --- Sheet: CUL (30 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 200
$ws.Range("K7").Value = 600
$ws.Range("M7").Value = -488
$ws.Range("H12").Value = 30328.705
$ws.Range("I12").Value = 49
$ws.Range("J12").Value = 103000
$ws.Range("K12").Value = 147
$ws.Range("L12").Value = 309000
$ws.Range("M12").Value = 26
$ws.Range("N12").Value = -309346
$ws.Range("H64").Value = 8366.538
$ws.Range("J64").Value = 9387.727999999999
$ws.Range("L64").Value = 28163.184
$ws.Range("N64").Value = -28703.184
$ws.Range("H67").Value = 8366.538
$ws.Range("J67").Value = 9387.727999999999
$ws.Range("L67").Value = 28163.184
$ws.Range("N67").Value = -30035.184
$ws.Range("H122").Value = 4764328.5
$ws.Range("I122").Value = 16667719
$ws.Range("K122").Value = 150009471
$ws.Range("M122").Value = -150007021
$ws.Range("H137").Value = 5945.3687
$ws.Range("I137").Value = 4497.1665
$ws.Range("J137").Value = 6613.769
$ws.Range("K137").Value = 13491.4995
$ws.Range("L137").Value = 19841.307
$ws.Range("M137").Value = -8391.499500000002
$ws.Range("N137").Value = -30041.307

# --- Sheet: GSM (36 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 17919700
$ws.Range("I80").Value = 68519.71000000001
$ws.Range("J80").Value = 71473240
$ws.Range("K80").Value = 68519.71000000001
$ws.Range("L80").Value = 71473240
$ws.Range("M80").Value = -67521.71000000001
$ws.Range("N80").Value = -71475236
$ws.Range("H83").Value = 17919700
$ws.Range("I83").Value = 68519.71000000001
$ws.Range("J83").Value = 71473240
$ws.Range("K83").Value = 342598.55
$ws.Range("L83").Value = 357366200
$ws.Range("M83").Value = -337606.55
$ws.Range("N83").Value = -357376184
$ws.Range("H126").Value = 1218
$ws.Range("I126").Value = 1313.2222
$ws.Range("J126").Value = 1003.75
$ws.Range("K126").Value = 3939.6666
$ws.Range("L126").Value = 3011.25
$ws.Range("M126").Value = -1469.6666
$ws.Range("N126").Value = -7951.25
$ws.Range("H132").Value = 13785394
$ws.Range("I132").Value = 3672.8462
$ws.Range("J132").Value = 39380016
$ws.Range("K132").Value = 11018.5386
$ws.Range("L132").Value = 118140048
$ws.Range("M132").Value = -8488.5386
$ws.Range("N132").Value = -118145108
$ws.Range("H135").Value = 103229.164
$ws.Range("J135").Value = 103229.164
$ws.Range("L135").Value = 103229.164
$ws.Range("N135").Value = -113369.164
$ws.Range("H137").Value = 127495
$ws.Range("J137").Value = 127495
$ws.Range("L137").Value = 127495
$ws.Range("N137").Value = -137695

# --- Sheet: LTW (53 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8450.423000000001
$ws.Range("I7").Value = 3922.25
$ws.Range("K7").Value = 3922.25
$ws.Range("M7").Value = -3810.25
$ws.Range("H22").Value = 5974.4
$ws.Range("J22").Value = 8136.5713
$ws.Range("L22").Value = 8136.5713
$ws.Range("N22").Value = -8726.5713
$ws.Range("H27").Value = 5974.4
$ws.Range("J27").Value = 8136.5713
$ws.Range("L27").Value = 8136.5713
$ws.Range("N27").Value = -8350.5713
$ws.Range("H40").Value = 3989.0557
$ws.Range("I40").Value = 2445.5557
$ws.Range("J40").Value = 5532.5557
$ws.Range("K40").Value = 2445.5557
$ws.Range("L40").Value = 5532.5557
$ws.Range("M40").Value = -2309.5557
$ws.Range("N40").Value = -5804.5557
$ws.Range("H55").Value = 237.125
$ws.Range("J55").Value = 253.84616
$ws.Range("L55").Value = 253.84616
$ws.Range("N55").Value = -599.8461600000001
$ws.Range("H82").Value = 6058
$ws.Range("I82").Value = 2392
$ws.Range("J82").Value = 11190.4
$ws.Range("K82").Value = 2392
$ws.Range("L82").Value = 11190.4
$ws.Range("M82").Value = -2031
$ws.Range("N82").Value = -11912.4
$ws.Range("H85").Value = 6058
$ws.Range("I85").Value = 2392
$ws.Range("J85").Value = 11190.4
$ws.Range("K85").Value = 2392
$ws.Range("L85").Value = 11190.4
$ws.Range("M85").Value = -1144
$ws.Range("N85").Value = -13686.4
$ws.Range("H126").Value = 8450.423000000001
$ws.Range("I126").Value = 3922.25
$ws.Range("K126").Value = 11766.75
$ws.Range("M126").Value = -9296.75
$ws.Range("H133").Value = 84498.5
$ws.Range("J133").Value = 84498.5
$ws.Range("L133").Value = 84498.5
$ws.Range("N133").Value = -89558.5
$ws.Range("H136").Value = 2762.0908
$ws.Range("I136").Value = 1719.7333
$ws.Range("K136").Value = 5159.199900000001
$ws.Range("M136").Value = -2609.199900000001
$ws.Range("H140").Value = 94071
$ws.Range("J140").Value = 94071
$ws.Range("L140").Value = 94071
$ws.Range("N140").Value = -104431

# --- Sheet: WVR (19 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7180.0835
$ws.Range("I62").Value = 6894.4287
$ws.Range("K62").Value = 6894.4287
$ws.Range("M62").Value = -6270.4287
$ws.Range("H65").Value = 7180.0835
$ws.Range("I65").Value = 6894.4287
$ws.Range("K65").Value = 34472.14350000001
$ws.Range("M65").Value = -31352.14350000001
$ws.Range("H100").Value = 648.2917
$ws.Range("J100").Value = 849.8333
$ws.Range("L100").Value = 1699.6666
$ws.Range("N100").Value = -2781.6666
$ws.Range("H113").Value = 323.78946
$ws.Range("I113").Value = 343.41666
$ws.Range("J113").Value = 290.14285
$ws.Range("K113").Value = 1030.24998
$ws.Range("L113").Value = 870.4285500000001
$ws.Range("M113").Value = 1139.75002
$ws.Range("N113").Value = -5210.428550000001

